$d = $word.ActiveDocument

# 1. "To have reliable replication of the bug (Incorrect calculation of fine
#     for one day overdue book)"
#     -> "...of fine, fine is half the amount intended.)"
$d.Content.Find.Execute(
    "To have reliable replication of the bug (Incorrect calculation of fine for one day overdue book)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To have reliable replication of the bug (Incorrect calculation of fine, fine is half the amount intended.)",
    2) | Out-Null

# 2. Version table date "17/10/2020" -> "18/10/2020"
$d.Content.Find.Execute(
    "17/10/2020", $true, $false, $false, $false, $false, $true, 1, $false,
    "18/10/2020", 2) | Out-Null

# 3. Test script heading "1.1 return" -> "1.0 return"
$d.Content.Find.Execute(
    "1.1 return", $true, $false, $false, $false, $false, $true, 1, $false,
    "1.0 return", 2) | Out-Null

# 4. Teardown bullet replaced entirely with "Testing"
$d.Content.Find.Execute(
    "Answer the question regarding the damage of the book. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Testing", 2) | Out-Null

# 5. Remove the two screenshots near the "Results" heading / automated-test
#    section (the floating "Picture 2" anchor and the inline "Picture 3"),
#    leaving their paragraphs empty. The earlier inline picture further up
#    the document (the one illustrating the loan-date step) is untouched.
for ($i = $d.Shapes.Count; $i -ge 1; $i--) {
    $shp = $d.Shapes.Item($i)
    if ($shp.Name -eq "Picture 2") {
        $shp.Delete() | Out-Null
    }
}

# The screenshot that follows the "From automated test:" line is the one
# being dropped (the earlier "Result of the following should be" picture
# stays put).
$automatedTestEnd = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "From automated test:*") {
        $automatedTestEnd = $para.Range.End
        break
    }
}

if ($automatedTestEnd -ne $null) {
    for ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {
        $ils = $d.InlineShapes.Item($i)
        if ($ils.Range.Start -eq $automatedTestEnd) {
            $ils.Delete() | Out-Null
        }
    }
}
